# arbeitszeit.xlsx - October "18:00-24:00" shift entries added to week of 27.10-02.11
# (rows 41-44) plus selection bookkeeping. Mirrors the author's manual entry of hours
# worked on Mon/Tue/Wed/Thu of that week, which previously had placeholder 0h/"9:00-17:00".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Monday 27.10 - worked 18:30 - 24:00 (4.5h)
$ws.Range("D41").Value = 4.5
$ws.Range("E41").Value = "18:30 - 24:00"

# Tuesday 28.10 - worked 18:00-24:00 (5h)
$ws.Range("D42").Value = 5
$ws.Range("E42").Value = "18:00-24:00"

# Wednesday 29.10 - worked 18:00 - 24:00 (5h)
$ws.Range("D43").Value = 5
$ws.Range("E43").Value = "18:00 - 24:00"

# Thursday 30.10 - worked 18:00 - 24:00 (5h)
$ws.Range("D44").Value = 5
$ws.Range("E44").Value = "18:00 - 24:00"

# Restore cursor/selection to where the author left off editing
[void]$ws.Range("F55").Select()
